$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Set the "Runmode" column (E) values to "Yes" for every test row (2 through 33)
# so that multi browser test execution is enabled across the board.
$ws.Range("E2:E33").Value = "Yes"

# Update the active selection on the sheet to match the edit
$ws.Cells.Item(6, 4).Select() | Out-Null
